# Build a data frame that keeps only the log-transformed response-time
# columns: drop the raw "incongruent"/"congruent" columns (D, E) and
# promote the existing "congruent_log"/"incongruent_log" columns (F, G)
# into their place. Deleting columns D:E shifts F->D and G->E for every
# row automatically (including row 46, which has no incongruent_log
# value), and keeps the bold/centered header style intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D:E").Delete()
